# Auto-generated edit script updating cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.864.27"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "3.529.74"
$ws.Range("E3").Value = "  -0.98%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'613.88"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").Value = "'173.47"
$ws.Range("E6").Value = "  +0.71%  "

$ws.Range("D7").Value = "'0.612"
$ws.Range("E7").Value = "  -0.99%  "

$ws.Range("D8").Value = "3.523.61"
$ws.Range("E8").Value = "  -1.01%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").Value = "'0.198"
$ws.Range("E10").Value = "  +0.85%  "

$ws.Range("E11").Value = "  -0.35%  "

$ws.Range("D12").Value = "'0.590"
$ws.Range("E12").Value = "  +0.73%  "

$ws.Range("D13").Value = "'46.66"
$ws.Range("E13").Value = "  -0.12%  "

$ws.Range("D14").Value = "'0.0000277"
$ws.Range("E14").Value = "  -0.02%  "

$ws.Range("D15").Value = "4.101.31"
$ws.Range("E15").Value = "  -0.90%  "

$ws.Range("D16").Value = "'8.45"
$ws.Range("E16").Value = "  +0.67%  "

$ws.Range("D17").Value = "'613.01"
$ws.Range("E17").Value = "  -0.74%  "

$ws.Range("D18").Value = "3.531.73"
$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").Value = "70.900.03"
$ws.Range("E19").Value = "  +0.29%  "

$ws.Range("E20").Value = "  +1.42%  "

$ws.Range("D21").Value = "'17.80"
$ws.Range("E21").Value = "  +2.26%  "

$ws.Range("D22").Value = "'0.887"
$ws.Range("E22").Value = "  +0.43%  "

$ws.Range("D23").Value = "'8.99"
$ws.Range("E23").Value = "  -4.60%  "

$ws.Range("D24").Value = "'15.78"
$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").Value = "'98.17"
$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("D26").Value = "'3.79"
$ws.Range("E26").Value = "  -1.40%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("E28").Value = "  -0.44%  "

$ws.Range("D29").Value = "'33.96"
$ws.Range("E29").Value = "  +1.11%  "

$ws.Range("D30").Value = "'9.20"
$ws.Range("E30").Value = "  +1.57%  "

$ws.Range("E31").Value = "  -0.62%  "

$ws.Range("E32").Value = "  -3.86%  "

$ws.Range("D33").Value = "'1.31"
$ws.Range("E33").Value = "  +0.22%  "

$ws.Range("D34").Value = "'6.87"
$ws.Range("E34").Value = "  -1.12%  "

$ws.Range("D35").Value = "'611.96"
$ws.Range("E35").Value = "  +6.51%  "

$ws.Range("E36").Value = "  -0.63%  "

$ws.Range("D37").Value = "'10.87"
$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("D38").Value = "'3.53"
$ws.Range("E38").Value = "  -2.20%  "

$ws.Range("D39").Value = "'0.0474"
$ws.Range("E39").Value = "  +0.66%  "

$ws.Range("D40").Value = "'57.05"
$ws.Range("E40").Value = "  -0.73%  "

$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("E42").Value = "  +2.00%  "

$ws.Range("D43").Value = "0.0₃0744"
$ws.Range("E43").Value = "  +6.15%  "

$ws.Range("D44").Value = "3.375.30"
$ws.Range("E44").Value = "  -0.31%  "

$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'2.97"
$ws.Range("E45").Value = "  +0.17%  "

$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.314"
$ws.Range("E46").Value = "  -2.01%  "

$ws.Range("D47").Value = "'32.24"
$ws.Range("E47").Value = "  -2.38%  "

$ws.Range("E48").Value = "  -1.49%  "

$ws.Range("E49").Value = "  +0.71%  "

$ws.Range("D50").Value = "'133.76"
$ws.Range("E50").Value = "  +0.17%  "

